$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "Cho"
$ws.Range("B35").Value = "Kyunghyun"
$ws.Range("C35").Value = 6
$ws.Range("D35").Value = 2

$ws.Range("A36").Value = "Bahdanau"
$ws.Range("B36").Value = "Dzmitri"
$ws.Range("C36").Value = 6
$ws.Range("D36").Value = 2

$ws.Range("B37").Select()
